$wb = $excel.ActiveWorkbook

# --- Sheet "Đơn phụ phẫu 1" (3rd sheet): insert a new detail row before the
#     "Tổng" (total) row and update the total row's aggregated values. ---
$ws3 = $wb.Worksheets.Item(3)

# Insert a new row at position 11, pushing the existing "Tổng" row down to 12
$ws3.Rows.Item(11).Insert()

# Fill in the new detail row (row 11)
$ws3.Cells.Item(11, 1).Value = "HD-LUXURY"
$ws3.Cells.Item(11, 2).Value = 696
# Force the date-shaped text to stay as plain text rather than being
# auto-converted into a date serial number, matching the other rows above.
$ws3.Cells.Item(11, 3).NumberFormat = "@"
$ws3.Cells.Item(11, 3).Value = "08-28-2024"
$ws3.Cells.Item(11, 4).Value = "SÓC TRĂNG"
$ws3.Cells.Item(11, 5).Value = "nguyễn thanh tuyền"
$ws3.Cells.Item(11, 6).Value = "Cá nhân"
$ws3.Cells.Item(11, 7).Value = "Nâng mũi"
$ws3.Cells.Item(11, 8).Value = "Kha Như Huỳnh "
$ws3.Cells.Item(11, 9).Value = 100000

# Update the "Tổng" row, now shifted to row 12
$ws3.Cells.Item(12, 2).Value = 10
$ws3.Cells.Item(12, 9).Value = 600000

# --- Sheet "Lương" (5th sheet): update computed salary summary values. ---
$ws5 = $wb.Worksheets.Item(5)

$ws5.Cells.Item(22, 2).Value = 23.5
$ws5.Cells.Item(23, 2).Value = 822500
$ws5.Cells.Item(24, 2).Value = 4616071.428571428
$ws5.Cells.Item(29, 2).Value = 600000
$ws5.Cells.Item(31, 2).Value = -1000000
$ws5.Cells.Item(34, 2).Value = 5438571.428571428
$ws5.Cells.Item(35, 2).Value = 5438571.428571428
